# Generated PowerShell/Excel COM script to append rows 21-24 to the "Artfynd" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21
$ws.Range("A21").Value = 111895224
$ws.Range("B21").Value = 90682
$ws.Range("C21").Value = "Ovaliderad"
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 2059
$ws.Range("F21").Value = "Skrovlig taggsvamp"
$ws.Range("G21").Value = "Hydnellum scabrosum"
$ws.Range("H21").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = ""
$ws.Range("J21").NumberFormat = "@"
$ws.Range("J21").Value = ""
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = ""
$ws.Range("N21").NumberFormat = "@"
$ws.Range("N21").Value = ""
$ws.Range("P21").Value = "Västanvik, Vrm"
$ws.Range("Q21").Value = 383386.3213553141
$ws.Range("R21").Value = 6664494.102709929
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = "Värmland"
$ws.Range("U21").Value = "Torsby"
$ws.Range("V21").Value = "Värmland"
$ws.Range("W21").Value = "Fryksände"
$ws.Range("Y21").NumberFormat = "@"
$ws.Range("Y21").Value = "2023-08-31"
$ws.Range("Z21").NumberFormat = "@"
$ws.Range("Z21").Value = "00:00"
$ws.Range("AA21").NumberFormat = "@"
$ws.Range("AA21").Value = "2023-08-31"
$ws.Range("AB21").NumberFormat = "@"
$ws.Range("AB21").Value = "00:00"
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AF21").NumberFormat = "@"
$ws.Range("AF21").Value = ""
$ws.Range("AG21").Value = $false
$ws.Range("AT21").NumberFormat = "@"
$ws.Range("AT21").Value = ""
$ws.Range("AW21").Value = "August Oljeqvist"
$ws.Range("AX21").Value = "August Oljeqvist, Jonas Göransson, Daniel Hertz Wallin , Amanda Evensen, Denise Persson, Johanna Klauss, Åsa Röstell, Edvin Johansson , Olavi Niemelä"
$ws.Range("AY21").NumberFormat = "@"
$ws.Range("AY21").Value = ""

# Row 22
$ws.Range("A22").Value = 111895144
$ws.Range("B22").Value = 56414
$ws.Range("C22").Value = "Ovaliderad"
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 100049
$ws.Range("F22").Value = "Spillkråka"
$ws.Range("G22").Value = "Dryocopus martius"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "1"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = ""
$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = "lockläte, övriga läten"
$ws.Range("N22").NumberFormat = "@"
$ws.Range("N22").Value = ""
$ws.Range("P22").Value = "Västanvik, Vrm"
$ws.Range("Q22").Value = 383214.7932507099
$ws.Range("R22").Value = 6664538.975714988
$ws.Range("S22").Value = 25
$ws.Range("T22").Value = "Värmland"
$ws.Range("U22").Value = "Torsby"
$ws.Range("V22").Value = "Värmland"
$ws.Range("W22").Value = "Fryksände"
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("Y22").Value = "2023-08-31"
$ws.Range("Z22").NumberFormat = "@"
$ws.Range("Z22").Value = "00:00"
$ws.Range("AA22").NumberFormat = "@"
$ws.Range("AA22").Value = "2023-08-31"
$ws.Range("AB22").NumberFormat = "@"
$ws.Range("AB22").Value = "00:00"
$ws.Range("AD22").Value = $false
$ws.Range("AE22").Value = $false
$ws.Range("AG22").Value = $false
$ws.Range("AT22").NumberFormat = "@"
$ws.Range("AT22").Value = ""
$ws.Range("AW22").Value = "August Oljeqvist"
$ws.Range("AX22").Value = "August Oljeqvist, Jonas Göransson, Daniel Hertz Wallin , Amanda Evensen, Denise Persson, Johanna Klauss, Åsa Röstell, Edvin Johansson , Olavi Niemelä"
$ws.Range("AY22").NumberFormat = "@"
$ws.Range("AY22").Value = ""

# Row 23
$ws.Range("A23").Value = 111895157
$ws.Range("B23").Value = 90709
$ws.Range("C23").Value = "Ovaliderad"
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 5448
$ws.Range("F23").Value = "Svartvit taggsvamp"
$ws.Range("G23").Value = "Phellodon connatus"
$ws.Range("H23").Value = "(Schultz) nom.prov"
$ws.Range("I23").NumberFormat = "@"
$ws.Range("I23").Value = ""
$ws.Range("J23").NumberFormat = "@"
$ws.Range("J23").Value = ""
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = ""
$ws.Range("N23").NumberFormat = "@"
$ws.Range("N23").Value = ""
$ws.Range("P23").Value = "Västanvik, Vrm"
$ws.Range("Q23").Value = 383310.7440082335
$ws.Range("R23").Value = 6664459.574209161
$ws.Range("S23").Value = 10
$ws.Range("T23").Value = "Värmland"
$ws.Range("U23").Value = "Torsby"
$ws.Range("V23").Value = "Värmland"
$ws.Range("W23").Value = "Fryksände"
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("Y23").Value = "2023-08-31"
$ws.Range("Z23").NumberFormat = "@"
$ws.Range("Z23").Value = "00:00"
$ws.Range("AA23").NumberFormat = "@"
$ws.Range("AA23").Value = "2023-08-31"
$ws.Range("AB23").NumberFormat = "@"
$ws.Range("AB23").Value = "00:00"
$ws.Range("AD23").Value = $false
$ws.Range("AE23").Value = $false
$ws.Range("AF23").NumberFormat = "@"
$ws.Range("AF23").Value = ""
$ws.Range("AG23").Value = $false
$ws.Range("AT23").NumberFormat = "@"
$ws.Range("AT23").Value = ""
$ws.Range("AW23").Value = "August Oljeqvist"
$ws.Range("AX23").Value = "August Oljeqvist, Jonas Göransson, Daniel Hertz Wallin , Amanda Evensen, Denise Persson, Johanna Klauss, Åsa Röstell, Edvin Johansson , Olavi Niemelä"
$ws.Range("AY23").NumberFormat = "@"
$ws.Range("AY23").Value = ""

# Row 24
$ws.Range("A24").Value = 111895200
$ws.Range("B24").Value = 90689
$ws.Range("C24").Value = "Ovaliderad"
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 5966
$ws.Range("F24").Value = "Motaggsvamp"
$ws.Range("G24").Value = "Sarcodon squamosus"
$ws.Range("H24").Value = "(Schaeff.) Quél."
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = ""
$ws.Range("J24").NumberFormat = "@"
$ws.Range("J24").Value = ""
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = ""
$ws.Range("N24").NumberFormat = "@"
$ws.Range("N24").Value = ""
$ws.Range("P24").Value = "Västanvik, Vrm"
$ws.Range("Q24").Value = 383318.0931039054
$ws.Range("R24").Value = 6664423.409828701
$ws.Range("S24").Value = 10
$ws.Range("T24").Value = "Värmland"
$ws.Range("U24").Value = "Torsby"
$ws.Range("V24").Value = "Värmland"
$ws.Range("W24").Value = "Fryksände"
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("Y24").Value = "2023-08-31"
$ws.Range("Z24").NumberFormat = "@"
$ws.Range("Z24").Value = "00:00"
$ws.Range("AA24").NumberFormat = "@"
$ws.Range("AA24").Value = "2023-08-31"
$ws.Range("AB24").NumberFormat = "@"
$ws.Range("AB24").Value = "00:00"
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AF24").NumberFormat = "@"
$ws.Range("AF24").Value = ""
$ws.Range("AG24").Value = $false
$ws.Range("AT24").NumberFormat = "@"
$ws.Range("AT24").Value = ""
$ws.Range("AW24").Value = "August Oljeqvist"
$ws.Range("AX24").Value = "August Oljeqvist, Jonas Göransson, Daniel Hertz Wallin , Amanda Evensen, Denise Persson, Johanna Klauss, Åsa Röstell, Edvin Johansson , Olavi Niemelä"
$ws.Range("AY24").NumberFormat = "@"
$ws.Range("AY24").Value = ""

# Reset number formatting on the cells above back to the default style
# (keeps their text type while removing the explicit "@" text format)
$clearRefs = @("I21","J21","K21","N21","Y21","Z21","AA21","AB21","AF21","AT21","AY21","I22","K22","L22","N22","Y22","Z22","AA22","AB22","AT22","AY22","I23","J23","K23","N23","Y23","Z23","AA23","AB23","AF23","AT23","AY23","I24","J24","K24","N24","Y24","Z24","AA24","AB24","AF24","AT24","AY24")
foreach ($r in $clearRefs) {
  $ws.Range($r).ClearFormats()
}
